# Auto-generated Excel COM-interop script
# Applies market-price/profit recalculation updates to the Ravana_Profits workbook
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 391.5

$ws.Range("I99").Value = 391.5

$ws.Range("K99").Value = 1174.5

$ws.Range("M99").Value = 323.5

$ws.Range("H101").Value = 831.6

$ws.Range("I101").Value = 831.6

$ws.Range("K101").Value = 2494.8

$ws.Range("M101").Value = -872.8000000000002

$ws.Range("H112").Value = 1321.6552

$ws.Range("I112").Value = 536.36365

$ws.Range("J112").Value = 1801.5555

$ws.Range("K112").Value = 1609.09095

$ws.Range("L112").Value = 5404.666499999999

$ws.Range("M112").Value = -501.09095

$ws.Range("N112").Value = -7620.666499999999

$ws.Range("H115").Value = 5295.125

$ws.Range("I115").Value = 5295.125

$ws.Range("K115").Value = 15885.375

$ws.Range("M115").Value = -14318.375

$ws.Range("H116").Value = 4832.6665

$ws.Range("J116").Value = 4832.6665

$ws.Range("L116").Value = 4832.6665

$ws.Range("N116").Value = -11716.6665

$ws.Range("H127").Value = 958.2

$ws.Range("I127").Value = 768.5

$ws.Range("K127").Value = 2305.5

$ws.Range("M127").Value = 2654.5

$ws.Range("H129").Value = 1351.2222

$ws.Range("I129").Value = 1276.8334

$ws.Range("J129").Value = 1500

$ws.Range("K129").Value = 3830.5002

$ws.Range("L129").Value = 4500

$ws.Range("M129").Value = 1169.4998

$ws.Range("N129").Value = -14500

$ws.Range("H132").Value = 1260.5

$ws.Range("I132").Value = 1222.05

$ws.Range("J132").Value = 1645

$ws.Range("K132").Value = 3666.15

$ws.Range("L132").Value = 4935

$ws.Range("M132").Value = -1136.15

$ws.Range("N132").Value = -9995

$ws.Range("H135").Value = 1164.75

$ws.Range("I135").Value = 1164.8572

$ws.Range("J135").Value = 1164

$ws.Range("K135").Value = 10483.7148

$ws.Range("L135").Value = 10476

$ws.Range("M135").Value = -7948.7148

$ws.Range("N135").Value = -15546

$ws.Range("H137").Value = 4861.9375

$ws.Range("I137").Value = 2719.4

$ws.Range("K137").Value = 8158.200000000001

$ws.Range("M137").Value = -5608.200000000001

$ws.Range("H138").Value = 2609.709

$ws.Range("I138").Value = 1049.0312

$ws.Range("J138").Value = 4781.087

$ws.Range("K138").Value = 3147.0936

$ws.Range("L138").Value = 14343.261

$ws.Range("M138").Value = 1992.9064

$ws.Range("N138").Value = -24623.261

$ws.Range("H141").Value = 2597.0476

$ws.Range("I141").Value = 1882.7222

$ws.Range("K141").Value = 5648.1666

$ws.Range("M141").Value = -468.1665999999996

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4549.1875

$ws.Range("I32").Value = 3852.8333

$ws.Range("K32").Value = 3852.8333

$ws.Range("M32").Value = -3565.8333

$ws.Range("H122").Value = 1914.9166

$ws.Range("I122").Value = 1914.9166

$ws.Range("K122").Value = 5744.7498

$ws.Range("M122").Value = -3294.7498

$ws.Range("H132").Value = 2017.2439

$ws.Range("I132").Value = 1858.1052

$ws.Range("K132").Value = 5574.3156

$ws.Range("M132").Value = -3044.3156

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3544.875

$ws.Range("I134").Value = 3614.6

$ws.Range("K134").Value = 10843.8

$ws.Range("M134").Value = -8308.799999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1996

$ws.Range("I31").Value = 2118.111

$ws.Range("K31").Value = 2118.111

$ws.Range("M31").Value = -1823.111

$ws.Range("H34").Value = 1996

$ws.Range("I34").Value = 2118.111

$ws.Range("K34").Value = 2118.111

$ws.Range("M34").Value = -1916.111

$ws.Range("H58").Value = 2499.125

$ws.Range("I58").Value = 2332.5

$ws.Range("K58").Value = 2332.5

$ws.Range("M58").Value = -2129.5

$ws.Range("H94").Value = 0

$ws.Range("J94").Value = 0

$ws.Range("L94").Value = 0

$ws.Range("N94").ClearContents()

$ws.Range("H99").Value = 4328

$ws.Range("I99").Value = 4328

$ws.Range("J99").Value = 0

$ws.Range("K99").Value = 4328

$ws.Range("L99").Value = 0

$ws.Range("M99").Value = -2830

$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 1380.3334

$ws.Range("I107").Value = 1070.5

$ws.Range("K107").Value = 1070.5

$ws.Range("M107").Value = 849.5

$ws.Range("H126").Value = 4328

$ws.Range("I126").Value = 4328

$ws.Range("J126").Value = 0

$ws.Range("K126").Value = 12984

$ws.Range("L126").Value = 0

$ws.Range("M126").Value = -10514

$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3544.2222

$ws.Range("I132").Value = 3299.9285

$ws.Range("J132").Value = 4399.25

$ws.Range("K132").Value = 9899.7855

$ws.Range("L132").Value = 13197.75

$ws.Range("M132").Value = -7369.7855

$ws.Range("N132").Value = -18257.75

$ws.Range("H136").Value = 2499.125

$ws.Range("I136").Value = 2332.5

$ws.Range("K136").Value = 6997.5

$ws.Range("M136").Value = -4447.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 5168.8

$ws.Range("J98").Value = 4783

$ws.Range("L98").Value = 14349

$ws.Range("N98").Value = -17345

$ws.Range("H107").Value = 1423.875

$ws.Range("J107").Value = 198.71428

$ws.Range("L107").Value = 596.14284

$ws.Range("N107").Value = -4436.14284

$ws.Range("H129").Value = 2900

$ws.Range("J129").Value = 0

$ws.Range("L129").Value = 0

$ws.Range("N129").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4741.8887

$ws.Range("I80").Value = 4083.875

$ws.Range("K80").Value = 4083.875

$ws.Range("M80").Value = -3085.875

$ws.Range("H83").Value = 4741.8887

$ws.Range("I83").Value = 4083.875

$ws.Range("K83").Value = 20419.375

$ws.Range("M83").Value = -15427.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4810.625

$ws.Range("J40").Value = 4832.6665

$ws.Range("L40").Value = 4832.6665

$ws.Range("N40").Value = -5104.6665

$ws.Range("H46").Value = 0

$ws.Range("I46").Value = 0

$ws.Range("J46").Value = 0

$ws.Range("K46").Value = 0

$ws.Range("L46").Value = 0

$ws.Range("M46").ClearContents()

$ws.Range("N46").ClearContents()

$ws.Range("H100").Value = 1950

$ws.Range("I100").Value = 1950

$ws.Range("K100").Value = 1950

$ws.Range("M100").Value = -1409

$ws.Range("H132").Value = 2218.1667

$ws.Range("I132").Value = 1445.1

$ws.Range("K132").Value = 4335.299999999999

$ws.Range("M132").Value = -1805.299999999999

$ws.Range("H136").Value = 3319.5293

$ws.Range("I136").Value = 3186.6924

$ws.Range("K136").Value = 9560.0772

$ws.Range("M136").Value = -7010.0772

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2120.25

$ws.Range("I122").Value = 2160.3333

$ws.Range("J122").Value = 2000

$ws.Range("K122").Value = 6480.999899999999

$ws.Range("L122").Value = 6000

$ws.Range("M122").Value = -4030.999899999999

$ws.Range("N122").Value = -10900

$ws.Range("H126").Value = 2541.9

$ws.Range("I126").Value = 2565.875

$ws.Range("K126").Value = 7697.625

$ws.Range("M126").Value = -5227.625

$ws.Range("H132").Value = 1489.8695

$ws.Range("I132").Value = 1355.5

$ws.Range("J132").Value = 1973.6

$ws.Range("K132").Value = 4066.5

$ws.Range("M132").Value = -1536.5
